{"js": "// Fix various feedback issues: apply five targeted text corrections\n// in the resume's OOXML body, each located via an exact, case-sensitive\n// search so we only touch the intended run.\nconst body = context.document.body;\n\n// Edit 1\n{\n  const results = body.search(\"Experienced coding distributed applications and map-reduce programs in Java. I've also used Scala but am less familiar with it.\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Edit 1: source text not found\");\n  }\n  results.items[0].insertText(\"I am experienced in coding distributed applications and map-reduce programs in Java. I've also used Scala but am less familiar with it.\", Word.InsertLocation.replace);\n}\n\n// Edit 2\n{\n  const results = body.search(\"SCRUM and hygienic coding\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Edit 2: source text not found\");\n  }\n  results.items[0].insertText(\"Scrum and hygienic coding\", Word.InsertLocation.replace);\n}\n\n// Edit 3\n{\n  const results = body.search(\"Experienced in a culture of culture of code reviews, code style guides, etc.\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Edit 3: source text not found\");\n  }\n  results.items[0].insertText(\"Participated in our culture of code reviews, code style guides, etc.\", Word.InsertLocation.replace);\n}\n\n// Edit 4\n{\n  const results = body.search(\"Independent Research with Prof. Emin Gun Sirer - (2001- 2002) Researched and developed an anonymizing peer-to-peer overlay network based on dining cryptographer nets called \u2018Herbivore\u2019. White paper available:\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Edit 4: source text not found\");\n  }\n  results.items[0].insertText(\"Independent Research with Prof. Emin Gun Sirer - (2001- 2002) Researched and developed an anonymizing peer-to-peer overlay network based on dining cryptographer nets called \u2018Herbivore\u2019. White paper:\", Word.InsertLocation.replace);\n}\n\n// Edit 5\n{\n  const results = body.search(\"Cornell University Neurobiology Department - (2005) Software engineer for audio experimentation and education cross platform workbench program called \u201cKoe\u0301\u201d, providing functionality similar to professional synthesis applications, such as Reaktor.\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Edit 5: source text not found\");\n  }\n  results.items[0].insertText(\"Cornell University Neurobiology Department - (2005) Worked as a software engineer on an audio research and education workbench program called \u201cKoe\u0301\u201d, providing functionality similar to professional synthesis applications, such as Reaktor.\", Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "# Fix various feedback issues: apply five targeted text corrections\n# by locating each exact (case-sensitive) source phrase with Word's\n# Find object, then rewriting just that matched range's text so the\n# surrounding runs/paragraph formatting stay untouched.\n$d = $word.ActiveDocument\n\n# Edit 1\n$find = $d.Content.Find\n$find.Text = 'Experienced coding distributed applications and map-reduce programs in Java. I''ve also used Scala but am less familiar with it.'\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\nif (-not $find.Execute()) { throw \"Edit 1: source text not found\" }\n$find.Parent.Text = 'I am experienced in coding distributed applications and map-reduce programs in Java. I''ve also used Scala but am less familiar with it.'\n\n# Edit 2\n$find = $d.Content.Find\n$find.Text = 'SCRUM and hygienic coding'\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\nif (-not $find.Execute()) { throw \"Edit 2: source text not found\" }\n$find.Parent.Text = 'Scrum and hygienic coding'\n\n# Edit 3\n$find = $d.Content.Find\n$find.Text = 'Experienced in a culture of culture of code reviews, code style guides, etc.'\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\nif (-not $find.Execute()) { throw \"Edit 3: source text not found\" }\n$find.Parent.Text = 'Participated in our culture of code reviews, code style guides, etc.'\n\n# Edit 4\n$find = $d.Content.Find\n$find.Text = 'Independent Research with Prof. Emin Gun Sirer - (2001- 2002) Researched and developed an anonymizing peer-to-peer overlay network based on dining cryptographer nets called \u2018Herbivore\u2019. White paper available:'\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\nif (-not $find.Execute()) { throw \"Edit 4: source text not found\" }\n$find.Parent.Text = 'Independent Research with Prof. Emin Gun Sirer - (2001- 2002) Researched and developed an anonymizing peer-to-peer overlay network based on dining cryptographer nets called \u2018Herbivore\u2019. White paper:'\n\n# Edit 5\n$find = $d.Content.Find\n$find.Text = 'Cornell University Neurobiology Department - (2005) Software engineer for audio experimentation and education cross platform workbench program called \u201cKoe\u0301\u201d, providing functionality similar to professional synthesis applications, such as Reaktor.'\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\nif (-not $find.Execute()) { throw \"Edit 5: source text not found\" }\n$find.Parent.Text = 'Cornell University Neurobiology Department - (2005) Worked as a software engineer on an audio research and education workbench program called \u201cKoe\u0301\u201d, providing functionality similar to professional synthesis applications, such as Reaktor.'\n"}
